$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A6").Value = "fb279db3-c387-4522-b066-07b4f10020be.md"
$ws1.Range("B6").Value = "e2e\fb279db3-c387-4522-b066-07b4f10020be.md"
$ws1.Range("E6").Value = "In Translation"
$ws1.Range("F6").Value = "In Translation"
$ws1.Range("G6").Value = "2016-10-14 07:37:40"

$ws1.Range("A7").Value = "5685cf38-4c69-4098-a2ff-8993427d9e74.md"
$ws1.Range("B7").Value = "e2e\5685cf38-4c69-4098-a2ff-8993427d9e74.md"
$ws1.Range("E7").Value = "Ready for handoff"
$ws1.Range("F7").Value = "Ready for handoff"
$ws1.Range("G7").Value = "2016-10-14 07:31:18"

# --- Sheet "zh-cn" ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A6").Value = "fb279db3-c387-4522-b066-07b4f10020be.md"
$ws2.Range("C6").Value = "In Translation"
$ws2.Range("G6").Value = "fb279db3-c387-4522-b066-07b4f10020be.663652f1cd374ac1997775048163ddb96a477949.zh-cn.xlf"
$ws2.Range("H6").Value = "2016-10-14 07:37:29"

$ws2.Range("A7").Value = "5685cf38-4c69-4098-a2ff-8993427d9e74.md"
$ws2.Range("C7").Value = "Ready for handoff"
$ws2.Range("G7").Value = "5685cf38-4c69-4098-a2ff-8993427d9e74.cb39d65c2879de697610e859672b22aa6064c84e.zh-cn.xlf"
$ws2.Range("H7").Value = "2016-10-14 07:31:07"

# --- Sheet "de-de" ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A6").Value = "fb279db3-c387-4522-b066-07b4f10020be.md"
$ws3.Range("C6").Value = "In Translation"
$ws3.Range("G6").Value = "fb279db3-c387-4522-b066-07b4f10020be.663652f1cd374ac1997775048163ddb96a477949.de-de.xlf"
$ws3.Range("H6").Value = "2016-10-14 07:37:40"

$ws3.Range("A7").Value = "5685cf38-4c69-4098-a2ff-8993427d9e74.md"
$ws3.Range("C7").Value = "Ready for handoff"
$ws3.Range("G7").Value = "5685cf38-4c69-4098-a2ff-8993427d9e74.cb39d65c2879de697610e859672b22aa6064c84e.de-de.xlf"
$ws3.Range("H7").Value = "2016-10-14 07:31:18"
